$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Firm, Time, Lawyers Registered (as text, matching source format)
$data = @(
    @('K1 Chamber', '05s', '0'),
    @('Pedersoli', '59s', '0'),
    @('MAS Law', '05s', '1'),
    @('Dillon Eustace', '07s', '2'),
    @('Nurmansyah And Muzdalifah', '14s', '1'),
    @('KISCH IP', '09s', '1'),
    @('Arthur Cox', '02min 22s', '2'),
    @('MSP', '45s', '1'),
    @('Thompson Dorfman Sweatman', '17s', '1'),
    @('Dompatent', '10s', '1'),
    @('Cobalt', '01min 40s', '2'),
    @('ShinAndKim', '15s', '1'),
    @('Gitti And Partners', '18s', '1'),
    @('Dottir', '04s', '0'),
    @('VB Advocates', '05s', '1'),
    @('Helmsman', '07s', '1'),
    @('Ellisons Solicitors', '20s', '1'),
    @('Control Risks', '03min 15s', '2'),
    @('Sangra', '10s', '1'),
    @('Duncan Cotterill', '30s', '1'),
    @('Allens', '03min 37s', '1'),
    @('Greenberg Traurig', '20s', '3'),
    @('Guantao Law', '52min 08s', '0'),
    @('Oxera', '02min 09s', '2'),
    @('Aera', '04s', '1'),
    @('Esche', '07s', '1'),
    @('SIRIUS', '10s', '0'),
    @('Huiye Law', '09s', '1'),
    @('Herbert Smith Freehills Kramer', '04s', '0'),
    @('Aron Tadmor Levy', '15s', '1'),
    @('Cassels', '01min 04s', '1'),
    @('Al Tamimi', '06s', '1'),
    @('Madrona', '14s', '0'),
    @('Hakun Law', '10s', '1'),
    @('Walkers', '32s', '3'),
    @('Stikeman Elliott', '07s', '1'),
    @('Schoenherr', '01min 13s', '2'),
    @('White And Case', '03min 10s', '3'),
    @('White And Case', '02min 58s', '0'),
    @('Dentons', '04min 27s', '1'),
    @('Clark Wilson', '08s', '1'),
    @('TEMPLARS', '11s', '1'),
    @('Tahota Law', '31s', '1'),
    @('Tuca Zbarcea', '05s', '1'),
    @('HFW', '06s', '3'),
    @('Hannes Snellman', '05s', '1'),
    @('Hill Dickinson', '07s', '2'),
    @('Bennett Jones', '18s', '1'),
    @('Byrne Wallace', '12s', '1'),
    @('Ramdas And Wong', '18s', '1'),
    @('Banki Haddock Fiora', '09s', '1'),
    @('Ellex', '01min 01s', '2'),
    @('Horten', '18s', '1'),
    @('Macpherson Kelley', '13s', '1'),
    @('Zhongzi Law', '08s', '1'),
    @('Higgs And Johnson', '30s', '1'),
    @('GrandwayLaw', '01min 57s', '1'),
    @('Blandy And Blandy', '11s', '1'),
    @('Carnelutti Law', '12s', '1'),
    @('Winston And Strawn', '14s', '2'),
    @('Spencer West', '46s', '3'),
    @('Sherrards  ', '37s', '1'),
    @('NPP Legal', '05s', '1'),
    @('Roschier', '32s', '1'),
    @('Gornitzky And Co', '14s', '1'),
    @('Asafo And Co', '34s', '2'),
    @('DW Fox Tucker', '11s', '1'),
    @('Clemens Law', '06s', '1'),
    @('Davies Ward Phillips And Vineberg', '10min 01s', '0'),
    @('Wolf Theiss', '29s', '3'),
    @('HY Leung And Co', '09s', '1'),
    @('BWB LLP', '09s', '1'),
    @('Arnesen IP', '05s', '1'),
    @('Santamarina And Steta', '17s', '1'),
    @('Ogier', '36s', '3'),
    @('Dale And Lessmann', '06s', '1'),
    @('EBN', '46s', '1'),
    @('Bae Kim And Lee', '01min 02s', '1'),
    @('Blakes', '06s', '1'),
    @('Fangda Partners', '01min 40s', '2'),
    @('Crowell And Moring', '02min 06s', '3'),
    @('Onsagers', '23s', '1'),
)

# Scratch cell used to coerce numeric-looking strings (column C) into text
# cells without altering the existing cell style (s=2, General format).
$scratch = $ws.Cells.Item(1000, 26)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]

    $scratch.NumberFormat = '@'
    $scratch.Value = $data[$i][2]
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($r, 3).PasteSpecial(-4163) | Out-Null
}

$scratch.Delete() | Out-Null
$excel.CutCopyMode = 0
